# exp2_run1_PCR.xlsx - "updated results for exp2"
#
# The underlying change is a correction of a mislabeled target name in the
# "Cq results" sheet: rows 2-57 of column C were re-typed from "GABB45"
# (a typo) to the correct "GADD45". Re-typing the cells causes Excel to
# record a (visually identical) new cell style for those cells. The
# active sheet/selection was also left on "Cq results" (instead of
# "Cq summary"), scrolled to row 43, with C2:C57 selected. Finally,
# iterative calculation was turned on (1 iteration).

$wb = $excel.ActiveWorkbook

# --- Cq results sheet: fix the mislabeled target name -----------------
$ws = $wb.Worksheets.Item("Cq results")
$ws.Activate()

$range = $ws.Range("C2:C57")
$range.Value = "GADD45"
# Re-typing the text in real Excel produced a new (but visually
# identical) style record for these cells - nudge the font size to force
# the same kind of new style entry.
$range.Font.Size = 8

# Leave the view/selection the way it was left when the file was saved:
# scrolled down so row 43 is at the top, with C2:C57 selected.
$range.Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1

# --- Workbook calculation options --------------------------------------
# Turn on iterative calculation with a single iteration.
$excel.Iteration = $true
$excel.MaxIterations = 1
